# Updated symbol list (Price / Volume(1h) columns) for the cryptos sheet,
# matching the "Updated symbol list on Fri Jan 20 05:54:22 UTC 2023 with
# GitHub Actions" data refresh. Only column D (Price) and column E
# (Volume(1h)) values change; everything else on the sheet is untouched.
#
# NumberFormat is forced to "@" (Text) before each write so that Excel
# keeps these numeric-looking / percent-looking strings as literal text
# (matching the source inlineStr cells) instead of silently re-parsing
# them into floating point numbers or percentage-formatted numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "289.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.29%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.44%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.918"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.42%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07399"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.69%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.268"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "26.85%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.716"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.68%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.750"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.06%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9118"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.77%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08845"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "14.61%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1687"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.29%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08239"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.92%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.82%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09951"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.67%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001498"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.53%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005843"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.10%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.493"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.37%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.087"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.31%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.48%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.08%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.970"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.63%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.16%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04557"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.91%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001209"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.43%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004582"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "14.35%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.05%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003398"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-95.49%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01592"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.12%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04473"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.23%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007301"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.62%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009569"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "24.67%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1325"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.38%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002233"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.84%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008885"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.52%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006079"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.58%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.05%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.200"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-2.02%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002002"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-33.29%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.05%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.05%"
